$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header ---
$ws.Range("A1").Value = "IdolStatType"

# --- Existing enum values renamed (rows 2-4) ---
$ws.Range("A2").Value = "Cute"
$ws.Range("B2").Value = 1

$ws.Range("A3").Value = "Cool"
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = "Sexy"
$ws.Range("B4").Value = 3

# --- New enum values (rows 5-8) ---
$ws.Range("A5").Value = "beauty"
$ws.Range("B5").Value = 4

$ws.Range("A6").Value = "vocal"
$ws.Range("B6").Value = 5

$ws.Range("A7").Value = "dance"
$ws.Range("B7").Value = 6

$ws.Range("A8").Value = "humor"
$ws.Range("B8").Value = 7

# --- Row 9: "intelligent" with a styled last character ---
$ws.Range("A9").Value = "intelligent"
$ws.Range("B9").Value = 8

# Whole-cell font for A9 (Microsoft YaHei, gray)
$ws.Range("A9").Font.Color = 6710886
$ws.Range("A9").Font.Name = "Microsoft YaHei"

# Rich-text run on just the final "t" character (position 11, length 1)
$run = $ws.Range("A9").Characters(11, 1)
$run.Font.Color = 6710886
$run.Font.Name = "맑은 고딕"

# --- Column width for column A ---
$ws.Columns("A").ColumnWidth = 11.43

# --- Selection ---
[void]$ws.Range("D7").Select()
